# cleanup both tts and stt to not use pipecat
# The "processing_time" column (E) is no longer populated for these rows;
# "ttfb" (D) gets a refreshed measurement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# elevenlabs row
$ws.Range("D2").Value = 1.362950086593628
$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"

# openai row
$ws.Range("D3").Value = 1.769400477409363
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"
